$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes all existing rows down by one)
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with the latest date entry (force text so the
# date-like string isn't auto-converted to a date serial number)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-13"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
